$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the orphan row 13 (a label-less row whose B/C cells only held the docente
# identification string). This shifts every row below it up by one, which is what
# produces the new row heights / dimension (A1:C21) seen in the target sheet.
$ws.Rows.Item(13).Delete()

# Row 10 (Objetivos:) now carries the docente identification text instead of the
# long "Apresentar ao aluno..." objectives paragraph.
$ws.Range("B10").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("C10").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

# Row 13 (Programa resumido:, formerly row 14) now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (Programa:, formerly row 16) now holds the activation date string. Force the
# cell to text first so Excel does not reinterpret the dd/mm/yyyy text as a date
# serial, then copy the plain (General, unwrapped-by-numberformat) look of an
# existing text cell back over it so the cell keeps its original column style.
$ws.Range("B15").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("C15").Value = "01/01/2012"
$ws.Range("B10").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 18 (Método:, formerly row 19) now holds the docente identification text.
$ws.Range("B18").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("C18").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

# Row 19 (Critério:, formerly row 20) now holds the teaching-method text.
$ws.Range("B19").Value = "Aulas expositivas; microcomputadores; seminários; visitas técnicas."
$ws.Range("C19").Value = "Aulas expositivas; microcomputadores; seminários; visitas técnicas."

# Row 20 (Norma de recuperação:, formerly row 21) now holds the grading-criteria text.
$ws.Range("B20").Value = "Média ponderada de notas de provas e seminários."
$ws.Range("C20").Value = "Média ponderada de notas de provas e seminários."

# Row 21 (Bibliografia:, formerly row 22) now holds the recovery-norm text instead of
# the bibliography (the bibliography text is dropped entirely, matching the shrunk
# shared-string table / the sheet shrinking to 21 rows).
$ws.Range("B21").Value = "Prova única com nota igual ou superior a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota igual ou superior a 5,0 (cinco)."
